# The commit swaps the contents of ppt/theme/theme1.xml (the theme used
# by the Slide Master -> every slide) and ppt/theme/theme2.xml (the theme
# used by the Notes Master): theme1 goes from the "Integral" colour
# scheme to the stock "Office Theme" colour scheme, and theme2 goes the
# other way. The font scheme (Arial-based "Office" fonts) and the format
# scheme (fills/lines/effects) are already byte-for-byte identical
# between the two theme parts, so the only real content difference is
# the 12-slot colour scheme (and the cosmetic name="" attributes, which
# aren't reachable from the exposed object model).
#
# The presentation's theme colours are exposed on PowerPoint's modern,
# 12-colour ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) via any Slide (they all resolve to the one Slide Master
# theme, i.e. ppt/theme/theme1.xml). Re-point every slot to the matching
# "Office Theme" RGB value to perform that half of the swap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (slot, target "Office Theme" RGB as a VBA RGB() long: R + G*256 + B*65536)
$officeThemeRgb = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeRgb[$i - 1]
}
